$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated S (TOTAL_SS_LOAD) and T (CONTESTABLE_ENERGY) values for rows 2-18,
# plus V (WESM_RATE) and W (CURRENT_RATE) for row 18, and T for rows 19-25.

$ws.Range("S2").Value = 78971
$ws.Range("T2").Value = 5623.7055

$ws.Range("S3").Value = 75703
$ws.Range("T3").Value = 5458.229

$ws.Range("S4").Value = 72637
$ws.Range("T4").Value = 5507.9745

$ws.Range("S5").Value = 70071
$ws.Range("T5").Value = 5475.498

$ws.Range("S6").Value = 70844
$ws.Range("T6").Value = 5469.52

$ws.Range("S7").Value = 73558
$ws.Range("T7").Value = 5575.703

$ws.Range("S8").Value = 73509
$ws.Range("T8").Value = 6141.8735

$ws.Range("S9").Value = 86371
$ws.Range("T9").Value = 7163.6355

$ws.Range("S10").Value = 103822
$ws.Range("T10").Value = 8647.460500000001

$ws.Range("S11").Value = 94388
$ws.Range("T11").Value = 12435.7905

$ws.Range("S12").Value = 97098
$ws.Range("T12").Value = 15471.022

$ws.Range("S13").Value = 97310
$ws.Range("T13").Value = 16384.8475

$ws.Range("S14").Value = 97015
$ws.Range("T14").Value = 16062.228

$ws.Range("S15").Value = 123772
$ws.Range("T15").Value = 16185.204

$ws.Range("S16").Value = 122973
$ws.Range("T16").Value = 16525.544

$ws.Range("S17").Value = 82093
$ws.Range("T17").Value = 16462.6175

$ws.Range("S18").Value = 56814
$ws.Range("T18").Value = 16761.7275
$ws.Range("V18").Value = 8856.747052777779
$ws.Range("W18").Value = 3.056410319552653

$ws.Range("T19").Value = 16532.1625
$ws.Range("T20").Value = 15890.3885
$ws.Range("T21").Value = 13683.775
$ws.Range("T22").Value = 11966.283
$ws.Range("T23").Value = 9398.494000000001
$ws.Range("T24").Value = 6574.9635
$ws.Range("T25").Value = 5831.934499999999
